$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: 2021年
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.7
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 102.3
$ws.Range("E7").Value = 101.2
$ws.Range("F7").Value = 102.2
$ws.Range("G7").Value = 102.8
$ws.Range("H7").Value = 101.6
$ws.Range("I7").Value = 101.3
$ws.Range("J7").Value = 101.3

# Row 8: 2022年 (only C8 populated)
$ws.Range("A8").Value = "2022年"
$ws.Range("C8").Value = 101.9

# Copy the formatting (style) used by the other year cells in column A
# onto the two new year-label cells, matching the existing s="1" cellXf.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
